# Conditional.xlsx — "Adding Module 04 Content"
#
# 1) Bump the Month Goal figure in C1.
# 2) Drop the bottom border that used to sit on the currency cells C6/C10
#    (they keep the same currency number format, just lose the border).
# 3) Add conditional formatting (green fill when a month's value meets/beats
#    the Month Goal in C1) to C5, C6 and C9:C10.
# 4) Restore the active selection to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Data change: Month Goal 1,200,000 -> 1,250,000 --------------------
$ws.Range("C1").Value = 1250000

# --- 2) Remove the bottom border from the currency totals in C6 and C10 ---
$ws.Range("C6").Borders.Item(9).LineStyle = -4142   # xlLineStyleNone
$ws.Range("C10").Borders.Item(9).LineStyle = -4142  # xlLineStyleNone

# --- 3) Conditional formatting: highlight months that hit the goal --------
# Green fill resolved from theme color "Accent3, Darker 25%" (the same
# swatch Excel itself renders for theme=6 / tint=-0.24994659260841701).
$goalFillColor = 3969911   # 0x77933C packed as BGR for Interior.Color

# Added in reverse row order so that, after re-priority below, the
# FormatConditions end up with the same priority numbering Excel produced
# (oldest/lowest-priority rule on top of the sheet, newest on the bottom-most
# range) while still appearing in the worksheet XML in row order.
$fcC9 = $ws.Range("C9:C10").FormatConditions.Add(2, 0, "=C9>=`$C`$1")
$fcC9.Interior.Color = $goalFillColor

$fcC6 = $ws.Range("C6").FormatConditions.Add(2, 0, "=C6>=`$C`$1")
$fcC6.Interior.Color = $goalFillColor

$fcC5 = $ws.Range("C5").FormatConditions.Add(2, 0, "=C5>=`$C`$1")
$fcC5.Interior.Color = $goalFillColor

$fcC5.Priority = 3
$fcC6.Priority = 2
$fcC9.Priority = 1

# --- 4) Move the active selection to E4 ------------------------------------
$ws.Range("E4").Select()
